$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.285.34"
$ws.Range("E2").Value = "  -2.64%  "
$ws.Range("D3").Value = "3.309.74"
$ws.Range("E3").Value = "  -2.80%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "558.03"
$ws.Range("E5").Value = "  -3.17%  "
$ws.Range("D6").Value = "142.54"
$ws.Range("E6").Value = "  -3.70%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").Value = "3.309.08"
$ws.Range("E8").Value = "  -2.86%  "
$ws.Range("E9").Value = "  -1.77%  "
$ws.Range("D10").Value = "7.85"
$ws.Range("E10").Value = "  -1.65%  "
$ws.Range("E11").Value = "  -3.10%  "
$ws.Range("D12").Value = "0.407"
$ws.Range("E12").Value = "  -1.30%  "
$ws.Range("D13").Value = "3.882.21"
$ws.Range("E13").Value = "  -2.69%  "
$ws.Range("E14").Value = "  +0.54%  "
$ws.Range("D15").Value = "26.98"
$ws.Range("E15").Value = "  -4.62%  "
$ws.Range("D16").Value = "3.311.67"
$ws.Range("E16").Value = "  -2.69%  "
$ws.Range("E17").Value = "  -2.66%  "
$ws.Range("D18").Value = "60.281.48"
$ws.Range("E18").Value = "  -2.73%  "
$ws.Range("E19").Value = "  -2.58%  "
$ws.Range("D20").Value = "14.43"
$ws.Range("E20").Value = "  +0.40%  "
$ws.Range("E21").Value = "  -2.82%  "
$ws.Range("D22").Value = "375.75"
$ws.Range("E22").Value = "  -1.07%  "
$ws.Range("D23").Value = "74.10"
$ws.Range("E23").Value = "  -0.46%  "
$ws.Range("E24").Value = "  -3.56%  "
$ws.Range("D26").Value = "3.441.58"
$ws.Range("E26").Value = "  -3.87%  "
$ws.Range("E27").Value = "  -7.19%  "
$ws.Range("D28").Value = "0.172"
$ws.Range("E28").Value = "  -4.44%  "
$ws.Range("E29").Value = "  -0.22%  "
$ws.Range("E30").Value = "  -4.80%  "
$ws.Range("E31").Value = "  -0.11%  "
$ws.Range("E32").Value = "  -2.76%  "
$ws.Range("E33").Value = "  -3.55%  "
$ws.Range("D34").Value = "22.60"
$ws.Range("E34").Value = "  -1.98%  "
$ws.Range("E35").Value = "  -4.10%  "
$ws.Range("E36").Value = "  -4.68%  "
$ws.Range("D37").Value = "166.63"
$ws.Range("E37").Value = "  -1.35%  "
$ws.Range("B38").Value = "Aptos"
$ws.Range("C38").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D38").Value = "6.76"
$ws.Range("E38").Value = "  -1.69%  "
$ws.Range("B39").Value = "ImmutableX"
$ws.Range("C39").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D39").Value = "1.53"
$ws.Range("E39").Value = "  -5.44%  "
$ws.Range("D40").Value = "26.88"
$ws.Range("E40").Value = "  -13.48%  "
$ws.Range("D41").Value = "3.340.21"
$ws.Range("E41").Value = "  -2.86%  "
$ws.Range("E42").Value = "  -5.12%  "
$ws.Range("D43").Value = "42.02"
$ws.Range("E43").Value = "  -0.72%  "
$ws.Range("E44").Value = "  -3.85%  "
$ws.Range("E45").Value = "  -3.44%  "
$ws.Range("E46").Value = "  -4.18%  "
$ws.Range("E47").Value = "  -3.72%  "
$ws.Range("D48").Value = "2.370.12"
$ws.Range("E48").Value = "  -6.82%  "
$ws.Range("E49").Value = "  +0.00%  "
$ws.Range("D50").Value = "6.54"
$ws.Range("E50").Value = "  -5.29%  "
$ws.Range("E51").Value = "  -4.98%  "
